# Commit: "added harvester and experiment design"
#
# The harvester column (B) was re-labelled from the placeholder
# "Retrofitted_0773" to the real harvester "S.GISH", and the
# experimentDesign column (D), previously blank, was populated with
# "90minuteInduction" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 22

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"   # D: experimentDesign
    $ws.Cells.Item($r, 2).Value = "S.GISH"               # B: harvester
}

# Restore the on-disk selection/scroll state captured in the saved file.
$ws.Range("C23:F28").Select() | Out-Null
